# Update "Forecast Comparison" sheet: insert a new "Week_Start_Date" column
# after "Week" (and before "ASIN"), shifting the rest of the columns right.
# Also normalize the "Week" labels from W01..W16 to W1..W16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column at B - this shifts ASIN..is_holiday_week (B..I) to (C..J)
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week-start dates (as text, matching the source file's inlineStr cells)
$weekStarts = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

# Force the new column to be treated as text so date-like strings
# ("2025-01-05") are not auto-converted into date serial numbers.
$ws.Columns.Item(2).NumberFormat = "@"

for ($i = 0; $i -lt $weekStarts.Length; $i++) {
    $row = $i + 2
    # Normalize "W01".."W16" -> "W1".."W16"
    $ws.Cells.Item($row, 1).Value = "W" + ($i + 1)
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $weekStarts[$i]
}
